$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coordinate values
$ws.Range("A2").Value = 1380
$ws.Range("B2").Value = 1000
$ws.Range("A4").Value = 1000
$ws.Range("A6").Value = 1320
$ws.Range("A13").Value = 630
$ws.Range("B13").Value = 1500

# Row 13 height change
$ws.Rows("13").RowHeight = 718.5

# Move view / selection to A4
$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
